# LoginList.xlsx edit: add "education" (school/degree/study/fromdate/ToDate/Program Desc)
# test data to the QaTeam3 sheet - rename the "Todate" column header to "ToDate" and
# store the from/to dates as literal text (Text-formatted cells) instead of numbers,
# then make QaTeam3 the active/selected sheet.

# G (fromdate) / H (ToDate) values, re-entered as text (dd/mm/yyyy-style strings)
# instead of the old bare numbers (11122021, 12122022, ...). Rows 3-10 are filled
# first (each introduces a fresh from/to date pair), then the header is renamed,
# and row 2 is filled in last (it reuses dates already used by rows 5 and 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QaTeam3")

$rowsInOrder = @(3, 4, 5, 6, 7, 8, 9, 10, 2)
$fromDateByRow = @{ 3 = "11/01/2002"; 4 = "11/01/2003"; 5 = "11/01/2004"; 6 = "11/01/2005"; 7 = "11/01/2006"; 8 = "11/01/2007"; 9 = "11/01/2008"; 10 = "11/01/2009"; 2 = "11/01/2004" }
$toDateByRow   = @{ 3 = "12/12/2005"; 4 = "12/12/2006"; 5 = "12/12/2007"; 6 = "12/12/2008"; 7 = "12/12/2009"; 8 = "12/12/2010"; 9 = "12/12/2011"; 10 = "12/12/2012"; 2 = "11/01/2005" }

foreach ($row in $rowsInOrder) {
    if ($row -eq 2) {
        # Rename the H column header from "Todate" to "ToDate" right before the
        # last row is (re)entered, matching the original authoring order.
        $ws.Cells.Item(1, 8).Value = "ToDate"
    }

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $fromDateByRow[$row]

    $hCell = $ws.Cells.Item($row, 8)
    $hCell.NumberFormat = "@"
    $hCell.Value = $toDateByRow[$row]
}

# Widen the G/H columns slightly now that they hold text dates.
$ws.Columns.Item(7).ColumnWidth = 11.86
$ws.Columns.Item(8).ColumnWidth = 11.86

# Make QaTeam3 the active sheet/tab, with H1 selected.
$ws.Activate()
$ws.Range("H1").Select()
